$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column A is formatted as text so years remain text values (matches t="str" in target)
$ws.Range("A1:A34").NumberFormat = "@"

$years = @(
    "2024",
    "2023",
    "2022",
    "2021",
    "2020",
    "2019",
    "2018",
    "2017",
    "2016",
    "2015",
    "2014",
    "2013",
    "2012",
    "2011",
    "2010",
    "2009",
    "2008",
    "2007",
    "2006",
    "2005",
    "2004",
    "2003",
    "2002",
    "2001",
    "2000",
    "1998",
    "1997",
    "1996",
    "1995",
    "1994",
    "1993",
    "1992",
    "1987"
)

$counts = @(
    76,
    68,
    53,
    52,
    51,
    47,
    37,
    35,
    30,
    21,
    20,
    18,
    19,
    16,
    17,
    10,
    10,
    9,
    3,
    7,
    6,
    3,
    2,
    5,
    3,
    1,
    2,
    2,
    4,
    4,
    1,
    2,
    1
)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $counts[$i]
}

